# Update the LR-pairs table (Gast -> Cckbr) with the newly computed TPM
# values. The sending-cluster set grew from {ECs, FAPs, Resolving-Mac} to
# {ECs, FAPs, MuSCs, Resolving-Mac}, and the target cluster now varies
# between ECs and FAPs for every sending cluster, so the table grows from
# 3 data rows to 8 data rows (A2:T9).
$data = @(
    ,@("ECs", "Gast", "Cckbr", "ECs", 3, 1, 3.513581333333333, 10.540744, 0.7556695049588564, 0.7556695049588565, 1, 0.3333333333333333, 0.1771103333333333, 0.531331, 0.6758517613986885, 0.6758517613986885, 0.6222915611404445, 5.600624050264, 0.5107205659617181, 0.5107205659617182)
    ,@("ECs", "Gast", "Cckbr", "FAPs", 3, 1, 3.513581333333333, 10.540744, 0.7556695049588564, 0.7556695049588565, 1, 0.3333333333333333, 0.08494466666666667, 0.254834, 0.3241482386013114, 0.3241482386013114, 0.2984599951662222, 2.686139956496, 0.2449489389971382, 0.2449489389971383)
    ,@("FAPs", "Gast", "Cckbr", "ECs", 1, 0.3333333333333333, 0.6282086666666666, 1.884626, 0.1351094758066973, 0.1351094758066973, 1, 0.3333333333333333, 0.1771103333333333, 0.531331, 0.6758517613986885, 0.6758517613986885, 0.1112622463562222, 1.001360217206, 0.09131397720560987, 0.09131397720560988)
    ,@("FAPs", "Gast", "Cckbr", "FAPs", 1, 0.3333333333333333, 0.6282086666666666, 1.884626, 0.1351094758066973, 0.1351094758066973, 1, 0.3333333333333333, 0.08494466666666667, 0.254834, 0.3241482386013114, 0.3241482386013114, 0.05336297578711111, 0.480266782084, 0.04379549860108743, 0.04379549860108744)
    ,@("MuSCs", "Gast", "Cckbr", "ECs", 1, 0.3333333333333333, 0.2157603333333333, 0.647281, 0.04640379396741573, 0.04640379396741574, 1, 0.3333333333333333, 0.1771103333333333, 0.531331, 0.6758517613986885, 0.6758517613986885, 0.03821338455677778, 0.343920461011, 0.03136208588845976, 0.03136208588845976)
    ,@("MuSCs", "Gast", "Cckbr", "FAPs", 1, 0.3333333333333333, 0.2157603333333333, 0.647281, 0.04640379396741573, 0.04640379396741574, 1, 0.3333333333333333, 0.08494466666666667, 0.254834, 0.3241482386013114, 0.3241482386013114, 0.01832768959488889, 0.164949206354, 0.01504170807895597, 0.01504170807895597)
    ,@("Resolving-Mac", "Gast", "Cckbr", "ECs", 1, 0.3333333333333333, 0.2920766666666667, 0.87623, 0.06281722526703037, 0.06281722526703037, 1, 0.3333333333333333, 0.1771103333333333, 0.531331, 0.6758517613986885, 0.6758517613986885, 0.05172979579222222, 0.46556816213, 0.04245513234290068, 0.04245513234290068)
    ,@("Resolving-Mac", "Gast", "Cckbr", "FAPs", 1, 0.3333333333333333, 0.2920766666666667, 0.87623, 0.06281722526703037, 0.06281722526703037, 1, 0.3333333333333333, 0.08494466666666667, 0.254834, 0.3241482386013114, 0.3241482386013114, 0.02481035509111111, 0.22329319582, 0.02036209292412969, 0.02036209292412969)
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
if (-not $ws) { $ws = $wb.Worksheets.Item(1) }

# Row 1 holds the column headers, so the first data row is row 2.
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $rowVals = $data[$i]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowVals[$c]
    }
}

Write-Output "done"
